$d = $word.ActiveDocument

# Each timecode run (bold, blue #4224E9) is immediately followed by a
# separate run containing exactly ten spaces with default formatting.
# The edit merges the timecode text and its trailing spaces into a single
# run: find the timecode + ten spaces as one match and replace it with
# itself, which collapses the pair into one run carrying the timecode
# run's bold/blue formatting (and marks the text xml:space="preserve").
$timecodes = @(
    "00:00:00.69",
    "00:00:02.50",
    "00:00:09.71",
    "00:00:12.09",
    "00:00:16.29",
    "00:00:20.15",
    "00:00:24.27",
    "00:00:26.93",
    "00:00:32.28",
    "00:00:37.94",
    "00:00:45.08",
    "00:00:49.10",
    "00:00:53.29",
    "00:01:00.18",
    "00:01:10.19",
    "00:01:16.01",
    "00:01:23.43",
    "00:01:29.33",
    "00:01:37.28",
    "00:01:51.97",
    "00:01:55.91",
    "00:02:03.44",
    "00:02:05.72",
    "00:02:10.98",
    "00:02:23.84",
    "00:02:27.35",
    "00:02:31.34",
    "00:02:33.06",
    "00:03:01.01",
    "00:03:04.06",
    "00:03:11.89",
    "00:03:14.72",
    "00:03:22.78",
    "00:03:30.62",
    "00:03:35.22",
    "00:03:42.63",
    "00:03:50.81",
    "00:03:55.47",
    "00:04:01.35",
    "00:04:10.08",
    "00:04:22.59",
    "00:04:25.00",
    "00:04:36.28",
    "00:04:42.34",
    "00:04:43.67"
)

foreach ($tc in $timecodes) {
    $needle = $tc + "          "
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
